# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the new columns, matching the style used by the other
# header cells in row 1 (bold text, thin border, centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the season record for every data row (rows 2-60).
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 30).Value = 100
    $ws.Cells.Item($r, 31).Value = 62
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "Updated dimension and added Wins/Losses/Ties columns"
